$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, matching the bold/border style used by the other headers (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# time_taken values for each data row (2..116)
$timeTaken = @{
    2 = "2021-10-05 13:39:15.542733"
    3 = "2021-10-05 13:39:15.542744"
    4 = "2021-10-05 13:39:15.542747"
    5 = "2021-10-05 13:39:15.542750"
    6 = "2021-10-05 13:39:15.542753"
    7 = "2021-10-05 13:39:15.542756"
    8 = "2021-10-05 13:39:15.542758"
    9 = "2021-10-05 13:39:15.542761"
    10 = "2021-10-05 13:39:15.542763"
    11 = "2021-10-05 13:39:15.542766"
    12 = "2021-10-05 13:39:15.542769"
    13 = "2021-10-05 13:39:15.542771"
    14 = "2021-10-05 13:39:15.542774"
    15 = "2021-10-05 13:39:15.542776"
    16 = "2021-10-05 13:39:15.542779"
    17 = "2021-10-05 13:39:15.542781"
    18 = "2021-10-05 13:39:15.542784"
    19 = "2021-10-05 13:39:15.542786"
    20 = "2021-10-05 13:39:15.542789"
    21 = "2021-10-05 13:39:15.542791"
    22 = "2021-10-05 13:39:15.542794"
    23 = "2021-10-05 13:39:15.542796"
    24 = "2021-10-05 13:39:15.542798"
    25 = "2021-10-05 13:39:15.542801"
    26 = "2021-10-05 13:39:15.542804"
    27 = "2021-10-05 13:39:15.542806"
    28 = "2021-10-05 13:39:15.542809"
    29 = "2021-10-05 13:39:15.542811"
    30 = "2021-10-05 13:39:15.542813"
    31 = "2021-10-05 13:39:15.542816"
    32 = "2021-10-05 13:39:15.542818"
    33 = "2021-10-05 13:39:15.542821"
    34 = "2021-10-05 13:39:15.542824"
    35 = "2021-10-05 13:39:15.542826"
    36 = "2021-10-05 13:39:15.542829"
    37 = "2021-10-05 13:39:15.542831"
    38 = "2021-10-05 13:39:15.542834"
    39 = "2021-10-05 13:39:15.542836"
    40 = "2021-10-05 13:39:15.542838"
    41 = "2021-10-05 13:39:15.542841"
    42 = "2021-10-05 13:39:15.542844"
    43 = "2021-10-05 13:39:15.542847"
    44 = "2021-10-05 13:39:15.542849"
    45 = "2021-10-05 13:39:15.542852"
    46 = "2021-10-05 13:39:15.542854"
    47 = "2021-10-05 13:39:15.542857"
    48 = "2021-10-05 13:39:15.542859"
    49 = "2021-10-05 13:39:15.542862"
    50 = "2021-10-05 13:39:15.542864"
    51 = "2021-10-05 13:39:15.542867"
    52 = "2021-10-05 13:39:15.542869"
    53 = "2021-10-05 13:39:15.542872"
    54 = "2021-10-05 13:39:15.542875"
    55 = "2021-10-05 13:39:15.542877"
    56 = "2021-10-05 13:39:15.542880"
    57 = "2021-10-05 13:39:15.542882"
    58 = "2021-10-05 13:39:15.542885"
    59 = "2021-10-05 13:39:15.542887"
    60 = "2021-10-05 13:39:15.542890"
    61 = "2021-10-05 13:39:15.542892"
    62 = "2021-10-05 13:39:15.542895"
    63 = "2021-10-05 13:39:15.542897"
    64 = "2021-10-05 13:39:15.542899"
    65 = "2021-10-05 13:39:15.542902"
    66 = "2021-10-05 13:39:15.542906"
    67 = "2021-10-05 13:39:15.542908"
    68 = "2021-10-05 13:39:15.542911"
    69 = "2021-10-05 13:39:15.542913"
    70 = "2021-10-05 13:39:15.542916"
    71 = "2021-10-05 13:39:15.542918"
    72 = "2021-10-05 13:39:15.542921"
    73 = "2021-10-05 13:39:15.542923"
    74 = "2021-10-05 13:39:15.542926"
    75 = "2021-10-05 13:39:15.542928"
    76 = "2021-10-05 13:39:15.542931"
    77 = "2021-10-05 13:39:15.542933"
    78 = "2021-10-05 13:39:15.542938"
    79 = "2021-10-05 13:39:15.542941"
    80 = "2021-10-05 13:39:15.542943"
    81 = "2021-10-05 13:39:15.542946"
    82 = "2021-10-05 13:39:15.542948"
    83 = "2021-10-05 13:39:15.542951"
    84 = "2021-10-05 13:39:15.542953"
    85 = "2021-10-05 13:39:15.542956"
    86 = "2021-10-05 13:39:15.542958"
    87 = "2021-10-05 13:39:15.542961"
    88 = "2021-10-05 13:39:15.542963"
    89 = "2021-10-05 13:39:15.542965"
    90 = "2021-10-05 13:39:15.542968"
    91 = "2021-10-05 13:39:15.542970"
    92 = "2021-10-05 13:39:15.542973"
    93 = "2021-10-05 13:39:15.542975"
    94 = "2021-10-05 13:39:15.542979"
    95 = "2021-10-05 13:39:15.542982"
    96 = "2021-10-05 13:39:15.542985"
    97 = "2021-10-05 13:39:15.542987"
    98 = "2021-10-05 13:39:15.542990"
    99 = "2021-10-05 13:39:15.542992"
    100 = "2021-10-05 13:39:15.542995"
    101 = "2021-10-05 13:39:15.542997"
    102 = "2021-10-05 13:39:15.543000"
    103 = "2021-10-05 13:39:15.543002"
    104 = "2021-10-05 13:39:15.543005"
    105 = "2021-10-05 13:39:15.543007"
    106 = "2021-10-05 13:39:15.543010"
    107 = "2021-10-05 13:39:15.543012"
    108 = "2021-10-05 13:39:15.543014"
    109 = "2021-10-05 13:39:15.543017"
    110 = "2021-10-05 13:39:15.543021"
    111 = "2021-10-05 13:39:15.543024"
    112 = "2021-10-05 13:39:15.543027"
    113 = "2021-10-05 13:39:15.543029"
    114 = "2021-10-05 13:39:15.543032"
    115 = "2021-10-05 13:39:15.543035"
    116 = "2021-10-05 13:39:15.543037"
}

foreach ($r in $timeTaken.Keys) {
    $ws.Cells.Item($r, 6).Value = $timeTaken[$r]
}

Write-Output "done"
